$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap the match data (columns F:V) between row 104 and row 105 ---
# Row 104 currently holds "Argentinos Jrs vs Independiente"; row 105 currently
# holds "Central Cordoba vs Tigre". The diff shows these two rows trade places
# (A:E - index/country/competition/season/date - stay put).

# Save row 104's current F:V values before overwriting anything.
$f104 = $ws.Range("F104").Value()
$g104 = $ws.Range("G104").Value()
$h104 = $ws.Range("H104").Value()
$i104 = $ws.Range("I104").Value()
$j104 = $ws.Range("J104").Value()
$k104 = $ws.Range("K104").Value()
$l104 = $ws.Range("L104").Value()
$m104 = $ws.Range("M104").Value()
$n104 = $ws.Range("N104").Value()
$o104 = $ws.Range("O104").Value()
$p104 = $ws.Range("P104").Value()
$q104 = $ws.Range("Q104").Value()
$r104 = $ws.Range("R104").Value()
$s104 = $ws.Range("S104").Value()
$t104 = $ws.Range("T104").Value()
$u104 = $ws.Range("U104").Value()
$v104 = $ws.Range("V104").Value()

# Save row 105's current F:V values too.
$f105 = $ws.Range("F105").Value()
$g105 = $ws.Range("G105").Value()
$h105 = $ws.Range("H105").Value()
$i105 = $ws.Range("I105").Value()
$j105 = $ws.Range("J105").Value()
$k105 = $ws.Range("K105").Value()
$l105 = $ws.Range("L105").Value()
$m105 = $ws.Range("M105").Value()
$n105 = $ws.Range("N105").Value()
$o105 = $ws.Range("O105").Value()
$p105 = $ws.Range("P105").Value()
$q105 = $ws.Range("Q105").Value()
$r105 = $ws.Range("R105").Value()
$s105 = $ws.Range("S105").Value()
$t105 = $ws.Range("T105").Value()
$u105 = $ws.Range("U105").Value()
$v105 = $ws.Range("V105").Value()

# Write row 105's old values into row 104.
$ws.Range("F104").Value = $f105
$ws.Range("G104").Value = $g105
$ws.Range("H104").Value = $h105
$ws.Range("I104").Value = $i105
$ws.Range("J104").Value = $j105
$ws.Range("K104").Value = $k105
$ws.Range("L104").Value = $l105
$ws.Range("M104").Value = $m105
$ws.Range("N104").Value = $n105
$ws.Range("O104").Value = $o105
$ws.Range("P104").Value = $p105
$ws.Range("Q104").Value = $q105
$ws.Range("R104").Value = $r105
$ws.Range("S104").Value = $s105
$ws.Range("T104").Value = $t105
$ws.Range("U104").Value = $u105
$ws.Range("V104").Value = $v105

# Write row 104's old values into row 105.
$ws.Range("F105").Value = $f104
$ws.Range("G105").Value = $g104
$ws.Range("H105").Value = $h104
$ws.Range("I105").Value = $i104
$ws.Range("J105").Value = $j104
$ws.Range("K105").Value = $k104
$ws.Range("L105").Value = $l104
$ws.Range("M105").Value = $m104
$ws.Range("N105").Value = $n104
$ws.Range("O105").Value = $o104
$ws.Range("P105").Value = $p104
$ws.Range("Q105").Value = $q104
$ws.Range("R105").Value = $r104
$ws.Range("S105").Value = $s104
$ws.Range("T105").Value = $t104
$ws.Range("U105").Value = $u104
$ws.Range("V105").Value = $v104

# --- Step 2: append the new match as row 135 ---
$ws.Range("A135").Value = 134
$ws.Range("B135").Value = "argentina"
$ws.Range("C135").Value = "copa-de-la-liga-profesional"
# "2023" looks numeric, so force it to stay text (matches the other rows,
# which all store the season as a text value) and then reset the format
# back to Normal so no stray number format is left applied to the cell.
$ws.Range("D135").NumberFormat = "@"
$ws.Range("D135").Value = "2023"
$ws.Range("D135").Style = "Normal"
$ws.Range("E135").Value = 45224.875
$ws.Range("F135").Value = "Arsenal Sarandi"
$ws.Range("G135").Value = 1
$ws.Range("H135").Value = "Colon Santa Fe"
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 2.96
$ws.Range("K135").Value = "20/10/2023 23:12"
$ws.Range("L135").Value = 3.11
$ws.Range("M135").Value = "25/10/2023 20:55"
$ws.Range("N135").Value = 3.01
$ws.Range("O135").Value = "20/10/2023 23:12"
$ws.Range("P135").Value = 3.02
$ws.Range("Q135").Value = "25/10/2023 20:55"
$ws.Range("R135").Value = 2.57
$ws.Range("S135").Value = "20/10/2023 23:12"
$ws.Range("T135").Value = 2.63
$ws.Range("U135").Value = "25/10/2023 20:55"
$ws.Range("V135").Value = "https://www.betexplorer.com/football/argentina/copa-de-la-liga-profesional/arsenal-sarandi-colon-santa-fe/8WUPnTUT/"

# Match the formatting used by the rest of the table: column A is bold /
# centered / bordered (same style as A134), column E is the datetime format
# used throughout column E (same style as E134).
$ws.Range("A134").Copy()
$ws.Range("A135").PasteSpecial(-4122)
$ws.Range("E134").Copy()
$ws.Range("E135").PasteSpecial(-4122)
